# Applies the "Updated symbol list" data refresh (Tue Jan 31 23:09:17 UTC 2023)
# to the cryptos worksheet: per-row Price/Volume(1h)/Hora updates, plus the
# FTXToken <-> KuCoinToken row swap (name/link/price/volume) on rows 6-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Plain text updates (coin name / link columns) -----------------------
$textUpdates = @{
    'B6' = 'KuCoinToken'
    'C6' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'B7' = 'FTXToken'
    'C7' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# --- Numeric-looking values that must stay stored as TEXT -----------------
# (Price / Volume(1h) / Hora columns hold strings like "311.51", "1.74%",
# "23" in the source workbook; force NumberFormat to Text first so Excel
# keeps writing them as text instead of silently coercing to numbers.)
$numericTextUpdates = @{
    'D2' = '311.51'
    'E2' = '1.74%'
    'G2' = '23'
    'D3' = '37.48'
    'E3' = '1.34%'
    'G3' = '23'
    'D4' = '5.120'
    'E4' = '0.78%'
    'G4' = '23'
    'D5' = '0.07885'
    'E5' = '2.15%'
    'G5' = '23'
    'D6' = '8.253'
    'E6' = '0.84%'
    'G6' = '23'
    'D7' = '1.898'
    'E7' = '0.96%'
    'G7' = '23'
    'D8' = '2.824'
    'E8' = '-11.03%'
    'G8' = '23'
    'D9' = '0.9221'
    'E9' = '-0.07%'
    'G9' = '23'
    'D10' = '0.1195'
    'E10' = '-2.04%'
    'G10' = '23'
    'D11' = '0.1923'
    'E11' = '3.47%'
    'G11' = '23'
    'D12' = '0.09054'
    'E12' = '3.47%'
    'G12' = '23'
    'D13' = '0.03315'
    'E13' = '-2.25%'
    'G13' = '23'
    'D14' = '0.09627'
    'E14' = '-0.84%'
    'G14' = '23'
    'D15' = '0.001380'
    'E15' = '0.89%'
    'G15' = '23'
    'D16' = '0.006035'
    'E16' = '1.21%'
    'G16' = '23'
    'D17' = '3.546'
    'E17' = '-1.56%'
    'G17' = '23'
    'D18' = '4.415'
    'E18' = '1.61%'
    'G18' = '23'
    'D19' = '0.3441'
    'E19' = '0.95%'
    'G19' = '23'
    'D20' = '5.233'
    'E20' = '4.24%'
    'G20' = '23'
    'D21' = '0.1283'
    'E21' = '1.70%'
    'G21' = '23'
    'D22' = '0.2588'
    'E22' = '-0.16%'
    'G22' = '23'
    'D23' = '0.04361'
    'E23' = '0.67%'
    'G23' = '23'
    'D24' = '0.001248'
    'E24' = '2.97%'
    'G24' = '23'
    'D25' = '0.004663'
    'E25' = '10.43%'
    'G25' = '23'
    'D26' = '0.0001360'
    'E26' = '0.51%'
    'G26' = '23'
    'D27' = '0.0003989'
    'E27' = '-98.11%'
    'G27' = '23'
    'G28' = '23'
    'G29' = '23'
    'G30' = '23'
    'G31' = '23'
    'G32' = '23'
    'G33' = '23'
    'G34' = '23'
    'G35' = '23'
    'G36' = '23'
    'G37' = '23'
    'G38' = '23'
    'D39' = '0.02251'
    'E39' = '3.99%'
    'G39' = '23'
    'D40' = '0.05064'
    'E40' = '3.65%'
    'G40' = '23'
    'D41' = '0.007450'
    'E41' = '-1.51%'
    'G41' = '23'
    'D42' = '0.009043'
    'E42' = '-8.25%'
    'G42' = '23'
    'D43' = '0.1354'
    'E43' = '1.24%'
    'G43' = '23'
    'D44' = '0.001950'
    'E44' = '-2.40%'
    'G44' = '23'
    'D45' = '0.009277'
    'E45' = '-6.34%'
    'G45' = '23'
    'D46' = '0.00006557'
    'E46' = '0.20%'
    'G46' = '23'
    'D47' = '0.00000000750'
    'E47' = '-0.23%'
    'G47' = '23'
    'D48' = '0.003317'
    'E48' = '10.37%'
    'G48' = '23'
    'D49' = '0.0009998'
    'E49' = '-23.26%'
    'G49' = '23'
    'D50' = '0.00002099'
    'E50' = '-0.23%'
    'G50' = '23'
    'D51' = '0.0001999'
    'E51' = '-0.23%'
    'G51' = '23'
}
foreach ($addr in $numericTextUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $numericTextUpdates[$addr]
}

